$wb = $excel.ActiveWorkbook
$wsRushing   = $wb.Worksheets.Item("Rushing")
$wsReceiving = $wb.Worksheets.Item("Receiving")

# ----------------------------------------------------------------------
# Sheet "Rushing" - Week 16 stat log + roster corrections
# Row count is unchanged (10 players), only values / a couple of names change
# ----------------------------------------------------------------------
$wsRushing.Range("D2").Value = 4

$wsRushing.Range("B3").Value = "R.Jones"
$wsRushing.Range("C3").Value = 139
$wsRushing.Range("D3").Value = 41
$wsRushing.Range("E3").Value = 20
$wsRushing.Range("F3").Value = 43

$wsRushing.Range("B4").Value = "G.Bernard"
$wsRushing.Range("C4").Value = 3
$wsRushing.Range("D4").Value = 0
$wsRushing.Range("E4").Value = 1
$wsRushing.Range("F4").Value = 0

$wsRushing.Range("B5").Value = "K.Vaughn"
$wsRushing.Range("C5").Value = 52
$wsRushing.Range("D5").Value = 18
$wsRushing.Range("E5").Value = 7
$wsRushing.Range("F5").Value = 12

$wsRushing.Range("B6").Value = "L.Bell"
$wsRushing.Range("C6").Value = 1
$wsRushing.Range("D6").Value = 1
$wsRushing.Range("E6").Value = 0
$wsRushing.Range("F6").Value = 0

$wsRushing.Range("B8").Value = "T.Johnson"
$wsRushing.Range("E8").Value = 1

$wsRushing.Range("B9").Value = "J.Darden"
$wsRushing.Range("C9").Value = 1
$wsRushing.Range("E9").Value = 0

$wsRushing.Range("B10").Value = "C.Grayson"

# Clear the stale selection box left on this sheet
$wsRushing.Range("A1").Select()

# ----------------------------------------------------------------------
# Sheet "Receiving" - remove players with no targets, add A.Brown,
# update Week 16 numbers, and drop the stray trailing row
# ----------------------------------------------------------------------

# Remove the trailing blank styled row (was row 16)
$wsReceiving.Rows(16).Delete()

# Remove C.Godwin's row (row 7) - released / no longer tracked
$wsReceiving.Rows(7).Delete()

# Remove L.Fournette's row (row 2) - released / no longer tracked
$wsReceiving.Rows(2).Delete()

# After the deletions above the remaining players shifted up to:
#  2 R.Jones, 3 G.Bernard, 4 K.Vaughn, 5 M.Evans, 6 S.Miller, 7 T.Johnson,
#  8 J.Darden, 9 C.Grayson, 10 B.Perriman, 11 R.Gronkowski, 12 O.Howard, 13 C.Brate

$wsReceiving.Range("A2").Value = 0
$wsReceiving.Range("A3").Value = 1
$wsReceiving.Range("A4").Value = 2
$wsReceiving.Range("A5").Value = 3
$wsReceiving.Range("A6").Value = 5
$wsReceiving.Range("A7").Value = 6
$wsReceiving.Range("A8").Value = 7
$wsReceiving.Range("A9").Value = 8
$wsReceiving.Range("A10").Value = 9
$wsReceiving.Range("A11").Value = 10
$wsReceiving.Range("A12").Value = 11
$wsReceiving.Range("A13").Value = 12

$wsReceiving.Range("B2").Value = "R.Jones"
$wsReceiving.Range("C2").Value = 33
$wsReceiving.Range("D2").Value = 30
$wsReceiving.Range("E2").Value = 2
$wsReceiving.Range("F2").Value = 1
$wsReceiving.Range("G2").Value = 11
$wsReceiving.Range("H2").Value = 9

$wsReceiving.Range("B3").Value = "G.Bernard"
$wsReceiving.Range("C3").Value = 26
$wsReceiving.Range("D3").Value = 21
$wsReceiving.Range("E3").Value = 2
$wsReceiving.Range("F3").Value = 1
$wsReceiving.Range("G3").Value = 7
$wsReceiving.Range("H3").Value = 4

$wsReceiving.Range("B4").Value = "K.Vaughn"
$wsReceiving.Range("C4").Value = 3
$wsReceiving.Range("D4").Value = 0
$wsReceiving.Range("E4").Value = 0
$wsReceiving.Range("F4").Value = 0
$wsReceiving.Range("G4").Value = 0
$wsReceiving.Range("H4").Value = 0

# Row 5 was M.Evans; he is replaced on the roster by A.Brown
$wsReceiving.Range("B5").Value = "A.Brown"
$wsReceiving.Range("C5").Value = 102
$wsReceiving.Range("D5").Value = 83
$wsReceiving.Range("E5").Value = 30
$wsReceiving.Range("F5").Value = 17
$wsReceiving.Range("G5").Value = 26
$wsReceiving.Range("H5").Value = 20

# Row 6 keeps its historical index value of 5 (gap in the numbering)
$wsReceiving.Range("B6").Value = "S.Miller"
$wsReceiving.Range("C6").Value = 20
$wsReceiving.Range("D6").Value = 15
$wsReceiving.Range("E6").Value = 6
$wsReceiving.Range("F6").Value = 4
$wsReceiving.Range("G6").Value = 2
$wsReceiving.Range("H6").Value = 1

$wsReceiving.Range("B7").Value = "T.Johnson"
$wsReceiving.Range("C7").Value = 35
$wsReceiving.Range("D7").Value = 23
$wsReceiving.Range("E7").Value = 6
$wsReceiving.Range("F7").Value = 3
$wsReceiving.Range("G7").Value = 5
$wsReceiving.Range("H7").Value = 2

$wsReceiving.Range("B8").Value = "J.Darden"
$wsReceiving.Range("C8").Value = 7
$wsReceiving.Range("D8").Value = 5
$wsReceiving.Range("E8").Value = 3
$wsReceiving.Range("F8").Value = 1
$wsReceiving.Range("G8").Value = 0
$wsReceiving.Range("H8").Value = 0

$wsReceiving.Range("B9").Value = "C.Grayson"
$wsReceiving.Range("C9").Value = 2
$wsReceiving.Range("D9").Value = 2
$wsReceiving.Range("E9").Value = 2
$wsReceiving.Range("F9").Value = 2
$wsReceiving.Range("G9").Value = 1
$wsReceiving.Range("H9").Value = 1

$wsReceiving.Range("B10").Value = "B.Perriman"
$wsReceiving.Range("C10").Value = 8
$wsReceiving.Range("D10").Value = 3
$wsReceiving.Range("E10").Value = 2
$wsReceiving.Range("F10").Value = 1
$wsReceiving.Range("G10").Value = 1
$wsReceiving.Range("H10").Value = 1

$wsReceiving.Range("B11").Value = "R.Gronkowski"
$wsReceiving.Range("C11").Value = 53
$wsReceiving.Range("D11").Value = 36
$wsReceiving.Range("E11").Value = 23
$wsReceiving.Range("F11").Value = 14
$wsReceiving.Range("G11").Value = 11
$wsReceiving.Range("H11").Value = 7

$wsReceiving.Range("B12").Value = "O.Howard"
$wsReceiving.Range("C12").Value = 18
$wsReceiving.Range("D12").Value = 13
$wsReceiving.Range("E12").Value = 3
$wsReceiving.Range("F12").Value = 1
$wsReceiving.Range("G12").Value = 2
$wsReceiving.Range("H12").Value = 1

$wsReceiving.Range("B13").Value = "C.Brate"
$wsReceiving.Range("C13").Value = 44
$wsReceiving.Range("D13").Value = 25
$wsReceiving.Range("E13").Value = 6
$wsReceiving.Range("F13").Value = 0
$wsReceiving.Range("G13").Value = 17
$wsReceiving.Range("H13").Value = 9

# Match the final active selection recorded on the Receiving sheet
$wsReceiving.Activate()
$wsReceiving.Range("F11").Select()
